# Update "想去人数" (want-to-go count) figures as refreshed by the
# gh-pages data generation script (commit 456a3b4).
#
# Sheet "展览" (Exhibitions)
#   F7  : 20    -> 21
#   F14 : 12290 -> 12292
#   F15 : 72    -> 74
#
# Sheet "演出" (Performances)
#   F2  : 112   -> 113
#
# Sheet "全部类型" (All types, aggregated view)
#   F4  : 112   -> 113
#   F9  : 20    -> 21
#   F16 : 12290 -> 12292
#   F18 : 72    -> 74

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F7").Value = 21
$wsExhibition.Range("F14").Value = 12292
$wsExhibition.Range("F15").Value = 74

$wsPerformance = $wb.Worksheets.Item("演出")
$wsPerformance.Range("F2").Value = 113

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 113
$wsAll.Range("F9").Value = 21
$wsAll.Range("F16").Value = 12292
$wsAll.Range("F18").Value = 74
